$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8346.727999999999
$ws.Range("I74").Value = 8281.4
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 8281.4
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -7345.4
$ws.Range("N74").Value = -10872

$ws.Range("H77").Value = 8346.727999999999
$ws.Range("I77").Value = 8281.4
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 41407
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -36727
$ws.Range("N77").Value = -54360

$ws.Range("H88").Value = 1445
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1445
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1445
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -2257

$ws.Range("H91").Value = 1445
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1445
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1445
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -4253

$ws.Range("H94").Value = 1585.7142
$ws.Range("I94").Value = 1100
$ws.Range("J94").Value = 1950
$ws.Range("K94").Value = 1100
$ws.Range("L94").Value = 1950
$ws.Range("M94").Value = -649
$ws.Range("N94").Value = -2852

$ws.Range("H101").Value = 2098.8
$ws.Range("I101").Value = 4498.5
$ws.Range("J101").Value = 499
$ws.Range("K101").Value = 13495.5
$ws.Range("L101").Value = 1497
$ws.Range("M101").Value = -11873.5
$ws.Range("N101").Value = -4741

$ws.Range("H137").Value = 1700
$ws.Range("I137").Value = 1725
$ws.Range("J137").Value = 1650
$ws.Range("K137").Value = 5175
$ws.Range("L137").Value = 4950
$ws.Range("M137").Value = -2625
$ws.Range("N137").Value = -10050

$ws.Range("H138").Value = 3255.2856
$ws.Range("I138").Value = 1748
$ws.Range("J138").Value = 3506.5
$ws.Range("K138").Value = 5244
$ws.Range("L138").Value = 10519.5
$ws.Range("M138").Value = -104
$ws.Range("N138").Value = -20799.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6000
$ws.Range("I2").Value = 6000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -5887
$ws.Range("N2").ClearContents()

$ws.Range("H97").Value = 511.3846
$ws.Range("I97").Value = 386.3
$ws.Range("J97").Value = 928.3333
$ws.Range("K97").Value = 386.3
$ws.Range("L97").Value = 928.3333
$ws.Range("M97").Value = 109.7
$ws.Range("N97").Value = -1920.3333

$ws.Range("H102").Value = 3400.1667
$ws.Range("I102").Value = 6045
$ws.Range("J102").Value = 2077.75
$ws.Range("K102").Value = 6045
$ws.Range("L102").Value = 2077.75
$ws.Range("M102").Value = -4423
$ws.Range("N102").Value = -5321.75

$ws.Range("H110").Value = 2141.8333
$ws.Range("I110").Value = 2349.7
$ws.Range("J110").Value = 1102.5
$ws.Range("K110").Value = 2349.7
$ws.Range("L110").Value = 1102.5
$ws.Range("M110").Value = -304.6999999999998
$ws.Range("N110").Value = -5192.5

$ws.Range("H116").Value = 6000
$ws.Range("I116").Value = 6000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -3706
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6000
$ws.Range("I3").Value = 6000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5886
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2680.2
$ws.Range("I31").Value = 2572.5715
$ws.Range("J31").Value = 2931.3333
$ws.Range("K31").Value = 2572.5715
$ws.Range("L31").Value = 2931.3333
$ws.Range("M31").Value = -2277.5715
$ws.Range("N31").Value = -3521.3333

$ws.Range("H34").Value = 2680.2
$ws.Range("I34").Value = 2572.5715
$ws.Range("J34").Value = 2931.3333
$ws.Range("K34").Value = 2572.5715
$ws.Range("L34").Value = 2931.3333
$ws.Range("M34").Value = -2370.5715
$ws.Range("N34").Value = -3335.3333

$ws.Range("H58").Value = 4879.8
$ws.Range("I58").Value = 2299.6667
$ws.Range("J58").Value = 8750
$ws.Range("K58").Value = 2299.6667
$ws.Range("L58").Value = 8750
$ws.Range("M58").Value = -2096.6667
$ws.Range("N58").Value = -9156

$ws.Range("H86").Value = 6591.7334
$ws.Range("I86").Value = 6284.143
$ws.Range("J86").Value = 6860.875
$ws.Range("K86").Value = 6284.143
$ws.Range("L86").Value = 6860.875
$ws.Range("M86").Value = -5161.143
$ws.Range("N86").Value = -9106.875

$ws.Range("H89").Value = 6591.7334
$ws.Range("I89").Value = 6284.143
$ws.Range("J89").Value = 6860.875
$ws.Range("K89").Value = 31420.715
$ws.Range("L89").Value = 34304.375
$ws.Range("M89").Value = -25804.715
$ws.Range("N89").Value = -45536.375

$ws.Range("H105").Value = 5384.5713
$ws.Range("I105").Value = 5449
$ws.Range("J105").Value = 4998
$ws.Range("K105").Value = 5449
$ws.Range("L105").Value = 4998
$ws.Range("M105").Value = -3702
$ws.Range("N105").Value = -8492

$ws.Range("H136").Value = 4879.8
$ws.Range("I136").Value = 2299.6667
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 6899.000100000001
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -4349.000100000001
$ws.Range("N136").Value = -31350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1415.8334
$ws.Range("I113").Value = 748
$ws.Range("J113").Value = 1749.75
$ws.Range("K113").Value = 2244
$ws.Range("L113").Value = 5249.25
$ws.Range("M113").Value = -74
$ws.Range("N113").Value = -9589.25

$ws.Range("H134").Value = 3663.3333
$ws.Range("I134").Value = 1995
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 5985
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -915
$ws.Range("N134").Value = -31140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 7999.5
$ws.Range("I97").Value = 7999.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 7999.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -7503.5
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 5054.375
$ws.Range("I102").Value = 5054.375
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5054.375
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3432.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 401.66666
$ws.Range("I16").Value = 405
$ws.Range("J16").Value = 395
$ws.Range("K16").Value = 405
$ws.Range("L16").Value = 395
$ws.Range("M16").Value = -235
$ws.Range("N16").Value = -735

$ws.Range("H127").Value = 100000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 100000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 499.15384
$ws.Range("I113").Value = 460.66666
$ws.Range("J113").Value = 585.75
$ws.Range("K113").Value = 1381.99998
$ws.Range("L113").Value = 1757.25
$ws.Range("M113").Value = 788.0000199999999
$ws.Range("N113").Value = -6097.25

$ws.Range("H122").Value = 2783.0625
$ws.Range("I122").Value = 2715.4
$ws.Range("J122").Value = 3798
$ws.Range("K122").Value = 8146.200000000001
$ws.Range("L122").Value = 11394
$ws.Range("M122").Value = -5696.200000000001
$ws.Range("N122").Value = -16294

$ws.Range("H126").Value = 1337.6
$ws.Range("I126").Value = 1337.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4012.8
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1542.8
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2114.4614
$ws.Range("I132").Value = 1636.5
$ws.Range("J132").Value = 2879.2
$ws.Range("K132").Value = 4909.5
$ws.Range("L132").Value = 8637.599999999999
$ws.Range("M132").Value = -2379.5
$ws.Range("N132").Value = -13697.6

$ws.Range("H136").Value = 7945.625
$ws.Range("I136").Value = 6094.3335
$ws.Range("J136").Value = 13499.5
$ws.Range("K136").Value = 18283.0005
$ws.Range("L136").Value = 40498.5
$ws.Range("M136").Value = -15733.0005
$ws.Range("N136").Value = -45598.5
